# Daily attendance processing - 2026-01-13 18:45:42
# For each data row, in column G ("Recorded By" style list of accounts),
# move the "System" entry (if present) to the front of the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "

        $hasSystem = $false
        foreach ($p in $parts) {
            if ($p.Equals("System")) { $hasSystem = $true }
        }

        if ($hasSystem) {
            $newParts = @("System")
            foreach ($p in $parts) {
                if (-not $p.Equals("System")) {
                    $newParts += $p
                }
            }
            $newVal = $newParts -join ", "
            if (-not $newVal.Equals($val)) {
                $cell.Value = $newVal
            }
        }
    }
}
